$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("First arithmetical Test")

# Row 20 header (COMMAND/param1/param2/result/exception): "exception" -> "disabled"
$ws.Range("E20").Value = "disabled"

# New DIV data block (row 21) inserted right after row 20's header, with a blank row 22 after it
$ws.Range("A21").Value = "DIV"
$ws.Range("B21").Value = 15
$ws.Range("C21").Value = 3
$ws.Range("D21").Value = 5
$ws.Range("E21").Formula = "=""true"""
$ws.Range("E21").Copy()
$ws.Range("E21").PasteSpecial(-4163)

# New header row (row 23): same look as the other "COMMAND/param1/param2/result/exception" headers
$ws.Range("A20:E20").Copy()
$ws.Range("A23:E23").PasteSpecial(-4122)
$ws.Range("A23").Value = "COMMAND"
$ws.Range("B23").Value = "param1"
$ws.Range("C23").Value = "param2"
$ws.Range("D23").Value = "result"
$ws.Range("E23").Value = "exception"

# Existing SUB block (previously row 21) shifted down to row 24
$ws.Range("A24").Value = "SUB"
$ws.Range("B24").Value = 3
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 0
$ws.Range("E24").Formula = "=""true"""
$ws.Range("E24").Copy()
$ws.Range("E24").PasteSpecial(-4163)

$ws.Range("F20").Select()
